$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.773585333333333
$ws.Range("H2").Value = 5.320756
$ws.Range("I2").Value = 0.3163248465685771
$ws.Range("J2").Value = 0.3163248465685771
$ws.Range("M2").Value = 1.294277666666667
$ws.Range("N2").Value = 3.882833
$ws.Range("O2").Value = 0.02669869520209519
$ws.Range("P2").Value = 0.02669869520209519
$ws.Range("Q2").Value = 2.295511886860889
$ws.Range("R2").Value = 20.659606981748
$ws.Range("S2").Value = 0.008445460663383968
$ws.Range("T2").Value = 0.008445460663383966
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.773585333333333
$ws.Range("H3").Value = 5.320756
$ws.Range("I3").Value = 0.3163248465685771
$ws.Range("J3").Value = 0.3163248465685771
$ws.Range("O3").Value = 0.9510146298089733
$ws.Range("P3").Value = 0.9510146298089732
$ws.Range("Q3").Value = 81.76674443378
$ws.Range("R3").Value = 735.90069990402
$ws.Range("S3").Value = 0.3008295568587956
$ws.Range("T3").Value = 0.3008295568587956
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.773585333333333
$ws.Range("H4").Value = 5.320756
$ws.Range("I4").Value = 0.3163248465685771
$ws.Range("J4").Value = 0.3163248465685771
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1528983333333333
$ws.Range("N4").Value = 0.458695
$ws.Range("O4").Value = 0.003154026453294555
$ws.Range("P4").Value = 0.003154026453294555
$ws.Range("Q4").Value = 0.2711782414911111
$ws.Range("R4").Value = 2.44060417342
$ws.Range("S4").Value = 0.0009976969339116335
$ws.Range("T4").Value = 0.0009976969339116333
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.773585333333333
$ws.Range("H5").Value = 5.320756
$ws.Range("I5").Value = 0.3163248465685771
$ws.Range("J5").Value = 0.3163248465685771
$ws.Range("M5").Value = 0.698547
$ws.Range("N5").Value = 2.095641
$ws.Range("O5").Value = 0.0144098085887325
$ws.Range("P5").Value = 0.0144098085887325
$ws.Range("Q5").Value = 1.238932713844
$ws.Range("R5").Value = 11.150394424596
$ws.Range("S5").Value = 0.004558180490913373
$ws.Range("T5").Value = 0.004558180490913372
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.773585333333333
$ws.Range("H6").Value = 5.320756
$ws.Range("I6").Value = 0.3163248465685771
$ws.Range("J6").Value = 0.3163248465685771
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.22895
$ws.Range("N6").Value = 0.68685
$ws.Range("O6").Value = 0.004722839946904511
$ws.Range("P6").Value = 0.004722839946904511
$ws.Range("Q6").Value = 0.4060623620666666
$ws.Range("R6").Value = 3.6545612586
$ws.Range("S6").Value = 0.001493951621572516
$ws.Range("T6").Value = 0.001493951621572516
$ws.Range("G7").Value = 1.759984333333334
$ws.Range("H7").Value = 5.279953000000001
$ws.Range("I7").Value = 0.3138990629553956
$ws.Range("J7").Value = 0.3138990629553955
$ws.Range("M7").Value = 1.294277666666667
$ws.Range("N7").Value = 3.882833
$ws.Range("O7").Value = 0.02669869520209519
$ws.Range("P7").Value = 0.02669869520209519
$ws.Range("Q7").Value = 2.277908416316556
$ws.Range("R7").Value = 20.501175746849
$ws.Range("S7").Value = 0.008380695406069397
$ws.Range("T7").Value = 0.008380695406069394
$ws.Range("G8").Value = 1.759984333333334
$ws.Range("H8").Value = 5.279953000000001
$ws.Range("I8").Value = 0.3138990629553956
$ws.Range("J8").Value = 0.3138990629553955
$ws.Range("O8").Value = 0.9510146298089733
$ws.Range("P8").Value = 0.9510146298089732
$ws.Range("Q8").Value = 81.13970412726502
$ws.Range("R8").Value = 730.2573371453851
$ws.Range("S8").Value = 0.2985226011539092
$ws.Range("T8").Value = 0.2985226011539091
$ws.Range("G9").Value = 1.759984333333334
$ws.Range("H9").Value = 5.279953000000001
$ws.Range("I9").Value = 0.3138990629553956
$ws.Range("J9").Value = 0.3138990629553955
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.1528983333333333
$ws.Range("N9").Value = 0.458695
$ws.Range("O9").Value = 0.003154026453294555
$ws.Range("P9").Value = 0.003154026453294555
$ws.Range("Q9").Value = 0.2690986712594445
$ws.Range("R9").Value = 2.421888041335
$ws.Range("S9").Value = 0.0009900459482256904
$ws.Range("T9").Value = 0.0009900459482256902
$ws.Range("G10").Value = 1.759984333333334
$ws.Range("H10").Value = 5.279953000000001
$ws.Range("I10").Value = 0.3138990629553956
$ws.Range("J10").Value = 0.3138990629553955
$ws.Range("M10").Value = 0.698547
$ws.Range("N10").Value = 2.095641
$ws.Range("O10").Value = 0.0144098085887325
$ws.Range("P10").Value = 0.0144098085887325
$ws.Range("Q10").Value = 1.229431776097
$ws.Range("R10").Value = 11.064885984873
$ws.Range("S10").Value = 0.004523225413369743
$ws.Range("T10").Value = 0.004523225413369741
$ws.Range("G11").Value = 1.759984333333334
$ws.Range("H11").Value = 5.279953000000001
$ws.Range("I11").Value = 0.3138990629553956
$ws.Range("J11").Value = 0.3138990629553955
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.22895
$ws.Range("N11").Value = 0.68685
$ws.Range("O11").Value = 0.004722839946904511
$ws.Range("P11").Value = 0.004722839946904511
$ws.Range("Q11").Value = 0.4029484131166667
$ws.Range("R11").Value = 3.62653571805
$ws.Range("S11").Value = 0.001482495033821636
$ws.Range("T11").Value = 0.001482495033821636
$ws.Range("G12").Value = 2.073278333333333
$ws.Range("H12").Value = 6.219835
$ws.Range("I12").Value = 0.3697760904760274
$ws.Range("J12").Value = 0.3697760904760274
$ws.Range("M12").Value = 1.294277666666667
$ws.Range("N12").Value = 3.882833
$ws.Range("O12").Value = 0.02669869520209519
$ws.Range("P12").Value = 0.02669869520209519
$ws.Range("Q12").Value = 2.683397843617222
$ws.Range("R12").Value = 24.150580592555
$ws.Range("S12").Value = 0.009872539132641829
$ws.Range("T12").Value = 0.009872539132641829
$ws.Range("G13").Value = 2.073278333333333
$ws.Range("H13").Value = 6.219835
$ws.Range("I13").Value = 0.3697760904760274
$ws.Range("J13").Value = 0.3697760904760274
$ws.Range("O13").Value = 0.9510146298089733
$ws.Range("P13").Value = 0.9510146298089732
$ws.Range("Q13").Value = 95.58334546167499
$ws.Range("R13").Value = 860.2501091550751
$ws.Range("S13").Value = 0.3516624717962686
$ws.Range("T13").Value = 0.3516624717962686
$ws.Range("G14").Value = 2.073278333333333
$ws.Range("H14").Value = 6.219835
$ws.Range("I14").Value = 0.3697760904760274
$ws.Range("J14").Value = 0.3697760904760274
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.1528983333333333
$ws.Range("N14").Value = 0.458695
$ws.Range("O14").Value = 0.003154026453294555
$ws.Range("P14").Value = 0.003154026453294555
$ws.Range("Q14").Value = 0.3170008017027777
$ws.Range("R14").Value = 2.853007215325
$ws.Range("S14").Value = 0.001166283571157231
$ws.Range("T14").Value = 0.001166283571157231
$ws.Range("G15").Value = 2.073278333333333
$ws.Range("H15").Value = 6.219835
$ws.Range("I15").Value = 0.3697760904760274
$ws.Range("J15").Value = 0.3697760904760274
$ws.Range("M15").Value = 0.698547
$ws.Range("N15").Value = 2.095641
$ws.Range("O15").Value = 0.0144098085887325
$ws.Range("P15").Value = 0.0144098085887325
$ws.Range("Q15").Value = 1.448282359915
$ws.Range("R15").Value = 13.034541239235
$ws.Range("S15").Value = 0.005328402684449386
$ws.Range("T15").Value = 0.005328402684449385
$ws.Range("G16").Value = 2.073278333333333
$ws.Range("H16").Value = 6.219835
$ws.Range("I16").Value = 0.3697760904760274
$ws.Range("J16").Value = 0.3697760904760274
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.22895
$ws.Range("N16").Value = 0.68685
$ws.Range("O16").Value = 0.004722839946904511
$ws.Range("P16").Value = 0.004722839946904511
$ws.Range("Q16").Value = 0.4746770744166666
$ws.Range("R16").Value = 4.272093669749999
$ws.Range("S16").Value = 0.001746393291510359
$ws.Range("T16").Value = 0.001746393291510359
